# Sheran Nancy TC.xlsx — add two new test-case rows (Home, Notification)
# to the "Sheet1" worksheet (rows 13 & 14), matching the upstream commit
# "Add files via upload".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$xlVAlignTop = -4160

# ---------------------------------------------------------------------
# Row 13 — TC_Home01 (reuses the pre-existing, already-formatted blank
# row 13, so only cell values/alignment need to be written)
# ---------------------------------------------------------------------
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "TC_Home01"
$ws.Range("C13").Value = "Home "
$ws.Range("D13").Value = "To check home option "
$ws.Range("E13").Value = "The user must check whether there is a home option"
$ws.Range("F13").Value = "Click Home option on the left corner`nof the page"
$ws.Range("G13").Value = "It shows lists of option such as my orders, my rewards etc."
$ws.Range("H13").Value = "Same as the result"
$ws.Range("I13").Value = "Pass"
$ws.Range("J13").Value = "1. Contains many options`n2. Each section are separated by a thin line `nto differentiate"

$ws.Range("A13:J13").VerticalAlignment = $xlVAlignTop
$ws.Range("F13").WrapText = $true
$ws.Range("J13").WrapText = $true

# ---------------------------------------------------------------------
# Row 14 — TC_Notifi01 (brand new row)
# ---------------------------------------------------------------------
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "TC_Notifi01"
$ws.Range("C14").Value = "Notification"
$ws.Range("D14").Value = "To check the notification section"
$ws.Range("E14").Value = "The user must check for  notification"
$ws.Range("F14").Value = "Click on the bell symbol to check for`nnotification"
$ws.Range("G14").Value = "Shows the notification for our account such as order confirm,`nOffers on products etc."
$ws.Range("H14").Value = "Same as the result"
$ws.Range("I14").Value = "Pass"
$ws.Range("J14").Value = "Notifies offers on our frequent purchase"

$ws.Range("A14:J14").VerticalAlignment = $xlVAlignTop
$ws.Range("D14").WrapText = $true
$ws.Range("E14").WrapText = $true
$ws.Range("F14").WrapText = $true
$ws.Range("G14").WrapText = $true
$ws.Range("J14").WrapText = $true

$ws.Rows.Item(14).RowHeight = 45

# ---------------------------------------------------------------------
# Move the view/selection the same way the author's session ended up:
# scrolled down a bit, with A15 selected.
# ---------------------------------------------------------------------
$null = $ws.Range("A15").Select()
